# Generate Report for Handback
# Applies the localization-status.xlsx "handback" update:
#  - Status text changes from "Ready for handoff" to "Handed back: in sync with en-US"
#    on the Overview sheet (zh-cn/de-de status columns) and on each language sheet's
#    Status cell.
#  - Each language sheet records the freshly generated handback artifacts: the
#    (hyperlinked) Latest Target File, the Latest Handback File, and the Latest
#    Handback DateTime.
#  - A few columns are widened to fit the newly-populated long file names.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$mdName     = "09c2564d-9e3f-4e64-b5f3-f06f832c5af7.md"
$mdUrl      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b14e9a00c0ced6aecc6b280ab6aa8ea9a4ac3748/e2e/09c2564d-9e3f-4e64-b5f3-f06f832c5af7.md"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText

$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("I2").Value = $mdName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl)
$wsZhCn.Range("J2").Value = "09c2564d-9e3f-4e64-b5f3-f06f832c5af7.200d9454f0c6994c646b0d2f93ae82418e5e23f9.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-05 03:07:00"

$wsZhCn.Columns.Item(3).ColumnWidth  = 29.166666666666668
$wsZhCn.Columns.Item(9).ColumnWidth  = 39.0
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("I2").Value = $mdName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl)
$wsDeDe.Range("J2").Value = "09c2564d-9e3f-4e64-b5f3-f06f832c5af7.200d9454f0c6994c646b0d2f93ae82418e5e23f9.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-05 03:07:14"

$wsDeDe.Columns.Item(3).ColumnWidth  = 29.166666666666668
$wsDeDe.Columns.Item(9).ColumnWidth  = 39.0
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664

Write-Output "Handback report generated"
